# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# c33b6286-7556-4d32-9fc9-7304f58620f1.md file is now "Ready for handoff"
# for both the zh-cn and de-de locales, including the new handoff
# timestamps and the (newly populated) error detail message.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ba6cdfbd7fd6d529511f07295fb4770c70cb479d/e2e/c33b6286-7556-4d32-9fc9-7304f58620f1.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f1419e4dce3c468c72927dd5904ac5a4bdae1903/e2e/c33b6286-7556-4d32-9fc9-7304f58620f1.md."

# --- Overview sheet: row 3 is the c33b6286 file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-21 22:56:47"

# --- zh-cn sheet: row 3 is the c33b6286 file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-21 22:56:43"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet: row 3 is the c33b6286 file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-21 22:56:47"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
